$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Car Names"
$ws.Range("J1").Value = "Predicted headform score (excluding blue points)"
